$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 715264.9
$ws.Range("I6").Value = 834416.5600000001
$ws.Range("K6").Value = 2503249.68
$ws.Range("M6").Value = -2503137.68
$ws.Range("H8").Value = 86.25
$ws.Range("I8").Value = 86.25
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 258.75
$ws.Range("L8").Value = 0
$ws.Range("M8").Value = -119.75
$ws.Range("N8").ClearContents()
$ws.Range("H39").Value = 636.7778
$ws.Range("I39").Value = 104.57143
$ws.Range("K39").Value = 313.71429
$ws.Range("M39").Value = -17.71429000000001
$ws.Range("H42").Value = 267.25
$ws.Range("J42").Value = 521.75
$ws.Range("L42").Value = 1565.25
$ws.Range("N42").Value = -2025.25
$ws.Range("H53").Value = 996.36365
$ws.Range("I53").Value = 283.5
$ws.Range("J53").Value = 1851.8
$ws.Range("K53").Value = 283.5
$ws.Range("L53").Value = 1851.8
$ws.Range("M53").Value = 353.5
$ws.Range("N53").Value = -3125.8
$ws.Range("H64").Value = 9428.286
$ws.Range("I64").Value = 8999.5
$ws.Range("K64").Value = 8999.5
$ws.Range("M64").Value = -8751.5
$ws.Range("H67").Value = 9428.286
$ws.Range("I67").Value = 8999.5
$ws.Range("K67").Value = 8999.5
$ws.Range("M67").Value = -8141.5
$ws.Range("H70").Value = 2000
$ws.Range("I70").Value = 2000
$ws.Range("J70").Value = 0
$ws.Range("K70").Value = 6000
$ws.Range("L70").Value = 0
$ws.Range("M70").Value = -5730
$ws.Range("N70").ClearContents()
$ws.Range("H73").Value = 2000
$ws.Range("I73").Value = 2000
$ws.Range("J73").Value = 0
$ws.Range("K73").Value = 6000
$ws.Range("L73").Value = 0
$ws.Range("M73").Value = -5064
$ws.Range("N73").ClearContents()
$ws.Range("H94").Value = 762
$ws.Range("I94").Value = 762
$ws.Range("K94").Value = 762
$ws.Range("M94").Value = -311
$ws.Range("H99").Value = 932.8
$ws.Range("J99").Value = 500
$ws.Range("L99").Value = 1500
$ws.Range("N99").Value = -4496
$ws.Range("H101").Value = 715.75
$ws.Range("J101").Value = 1166.3334
$ws.Range("L101").Value = 3499.0002
$ws.Range("N101").Value = -6743.0002
$ws.Range("H106").Value = 13338085
$ws.Range("I106").Value = 13338085
$ws.Range("K106").Value = 13338085
$ws.Range("M106").Value = -13337454
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 23281372
$ws.Range("I32").Value = 23245066
$ws.Range("K32").Value = 23245066
$ws.Range("M32").Value = -23244779
$ws.Range("H74").Value = 1629.3256
$ws.Range("I74").Value = 1538.8387
$ws.Range("J74").Value = 1863.0834
$ws.Range("K74").Value = 1538.8387
$ws.Range("L74").Value = 1863.0834
$ws.Range("M74").Value = -664.8387
$ws.Range("N74").Value = -3611.0834
$ws.Range("H77").Value = 1629.3256
$ws.Range("I77").Value = 1538.8387
$ws.Range("J77").Value = 1863.0834
$ws.Range("K77").Value = 7694.1935
$ws.Range("L77").Value = 9315.416999999999
$ws.Range("M77").Value = -3326.1935
$ws.Range("N77").Value = -18051.417
$ws.Range("H111").Value = 0
$ws.Range("J111").Value = 0
$ws.Range("L111").Value = 0
$ws.Range("N111").ClearContents()
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H38").Value = 98821.336
$ws.Range("J38").Value = 98821.336
$ws.Range("L38").Value = 98821.336
$ws.Range("N38").Value = -99653.336
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4651.5835
$ws.Range("I31").Value = 2234.7856
$ws.Range("K31").Value = 2234.7856
$ws.Range("M31").Value = -1939.7856
$ws.Range("H34").Value = 4651.5835
$ws.Range("I34").Value = 2234.7856
$ws.Range("K34").Value = 2234.7856
$ws.Range("M34").Value = -2032.7856
$ws.Range("H58").Value = 2781.9614
$ws.Range("I58").Value = 2577.9148
$ws.Range("K58").Value = 2577.9148
$ws.Range("M58").Value = -2374.9148
$ws.Range("H62").Value = 2182.5
$ws.Range("I62").Value = 2219.4
$ws.Range("J62").Value = 1998
$ws.Range("K62").Value = 2219.4
$ws.Range("L62").Value = 1998
$ws.Range("M62").Value = -1595.4
$ws.Range("N62").Value = -3246
$ws.Range("H65").Value = 2182.5
$ws.Range("I65").Value = 2219.4
$ws.Range("J65").Value = 1998
$ws.Range("K65").Value = 11097
$ws.Range("L65").Value = 9990
$ws.Range("M65").Value = -7977
$ws.Range("N65").Value = -16230
$ws.Range("H133").Value = 80326
$ws.Range("J133").Value = 80326
$ws.Range("L133").Value = 80326
$ws.Range("N133").Value = -85386
$ws.Range("H134").Value = 1779.2142
$ws.Range("I134").Value = 1437.2222
$ws.Range("J134").Value = 2394.8
$ws.Range("K134").Value = 4311.6666
$ws.Range("L134").Value = 7184.400000000001
$ws.Range("M134").Value = -1776.6666
$ws.Range("N134").Value = -12254.4
$ws.Range("H136").Value = 2781.9614
$ws.Range("I136").Value = 2577.9148
$ws.Range("K136").Value = 7733.7444
$ws.Range("M136").Value = -5183.7444
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H52").Value = 3056.1428
$ws.Range("J52").Value = 3056.1428
$ws.Range("L52").Value = 9168.428400000001
$ws.Range("N52").Value = -9700.428400000001
$ws.Range("H75").Value = 800
$ws.Range("J75").Value = 1000
$ws.Range("L75").Value = 3000
$ws.Range("N75").Value = -4996
$ws.Range("H78").Value = 800
$ws.Range("J78").Value = 1000
$ws.Range("L78").Value = 9000
$ws.Range("N78").Value = -18984
$ws.Range("H80").Value = 1999.5
$ws.Range("I80").Value = 1000
$ws.Range("K80").Value = 3000
$ws.Range("M80").Value = -2064
$ws.Range("H83").Value = 1999.5
$ws.Range("I83").Value = 1000
$ws.Range("K83").Value = 9000
$ws.Range("M83").Value = -4320
$ws.Range("H87").Value = 1499.5
$ws.Range("I87").Value = 1499.5
$ws.Range("K87").Value = 4498.5
$ws.Range("M87").Value = -3250.5
$ws.Range("H90").Value = 1499.5
$ws.Range("I90").Value = 1499.5
$ws.Range("K90").Value = 13495.5
$ws.Range("M90").Value = -7255.5
$ws.Range("H122").Value = 877.6
$ws.Range("J122").Value = 872
$ws.Range("L122").Value = 7848
$ws.Range("N122").Value = -12748
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 3131.389
$ws.Range("J126").Value = 3429.1667
$ws.Range("L126").Value = 10287.5001
$ws.Range("N126").Value = -15227.5001
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 11335.4
$ws.Range("I122").Value = 10231.19
$ws.Range("K122").Value = 30693.57
$ws.Range("M122").Value = -28243.57
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 1155354.6
$ws.Range("I4").Value = 1112877.8
$ws.Range("K4").Value = 1112877.8
$ws.Range("M4").Value = -1112764.8
$ws.Range("H103").Value = 59000
$ws.Range("J103").Value = 59000
$ws.Range("L103").Value = 59000
$ws.Range("N103").Value = -61344
$ws.Range("H107").Value = 790
$ws.Range("I107").Value = 737.5
$ws.Range("K107").Value = 2212.5
$ws.Range("M107").Value = -292.5
